$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column B (PO_Forecast) values for rows 2-84 per "4wk low sales check" adjustment
$ws.Cells.Item(2, 2).Value = 162
$ws.Cells.Item(3, 2).Value = 165
$ws.Cells.Item(4, 2).Value = 169
$ws.Cells.Item(5, 2).Value = 172
$ws.Cells.Item(6, 2).Value = 179
$ws.Cells.Item(7, 2).Value = 182
$ws.Cells.Item(8, 2).Value = 186
$ws.Cells.Item(9, 2).Value = 189
$ws.Cells.Item(10, 2).Value = 192
$ws.Cells.Item(11, 2).Value = 196
$ws.Cells.Item(12, 2).Value = 199
$ws.Cells.Item(13, 2).Value = 206
$ws.Cells.Item(14, 2).Value = 209
$ws.Cells.Item(15, 2).Value = 213
$ws.Cells.Item(16, 2).Value = 216
$ws.Cells.Item(17, 2).Value = 220
$ws.Cells.Item(18, 2).Value = 223
$ws.Cells.Item(19, 2).Value = 226
$ws.Cells.Item(20, 2).Value = 230
$ws.Cells.Item(21, 2).Value = 233
$ws.Cells.Item(22, 2).Value = 237
$ws.Cells.Item(23, 2).Value = 240
$ws.Cells.Item(24, 2).Value = 243
$ws.Cells.Item(25, 2).Value = 247
$ws.Cells.Item(26, 2).Value = 250
$ws.Cells.Item(27, 2).Value = 254
$ws.Cells.Item(28, 2).Value = 257
$ws.Cells.Item(29, 2).Value = 260
$ws.Cells.Item(30, 2).Value = 267
$ws.Cells.Item(31, 2).Value = 277
$ws.Cells.Item(32, 2).Value = 281
$ws.Cells.Item(33, 2).Value = 284
$ws.Cells.Item(34, 2).Value = 288
$ws.Cells.Item(35, 2).Value = 291
$ws.Cells.Item(36, 2).Value = 305
$ws.Cells.Item(37, 2).Value = 308
$ws.Cells.Item(38, 2).Value = 311
$ws.Cells.Item(39, 2).Value = 315
$ws.Cells.Item(40, 2).Value = 318
$ws.Cells.Item(41, 2).Value = 322
$ws.Cells.Item(42, 2).Value = 325
$ws.Cells.Item(43, 2).Value = 328
$ws.Cells.Item(44, 2).Value = 339
$ws.Cells.Item(45, 2).Value = 342
$ws.Cells.Item(46, 2).Value = 349
$ws.Cells.Item(47, 2).Value = 352
$ws.Cells.Item(48, 2).Value = 356
$ws.Cells.Item(49, 2).Value = 359
$ws.Cells.Item(50, 2).Value = 362
$ws.Cells.Item(51, 2).Value = 366
$ws.Cells.Item(52, 2).Value = 369
$ws.Cells.Item(53, 2).Value = 373
$ws.Cells.Item(54, 2).Value = 376
$ws.Cells.Item(55, 2).Value = 379
$ws.Cells.Item(56, 2).Value = 383
$ws.Cells.Item(57, 2).Value = 386
$ws.Cells.Item(58, 2).Value = 393
$ws.Cells.Item(59, 2).Value = 396
$ws.Cells.Item(60, 2).Value = 400
$ws.Cells.Item(61, 2).Value = 403
$ws.Cells.Item(62, 2).Value = 407
$ws.Cells.Item(63, 2).Value = 410
$ws.Cells.Item(64, 2).Value = 413
$ws.Cells.Item(65, 2).Value = 417
$ws.Cells.Item(66, 2).Value = 420
$ws.Cells.Item(67, 2).Value = 434
$ws.Cells.Item(68, 2).Value = 437
$ws.Cells.Item(69, 2).Value = 441
$ws.Cells.Item(70, 2).Value = 458
$ws.Cells.Item(71, 2).Value = 461
$ws.Cells.Item(72, 2).Value = 464
$ws.Cells.Item(73, 2).Value = 475
$ws.Cells.Item(74, 2).Value = 485
$ws.Cells.Item(75, 2).Value = 495
$ws.Cells.Item(76, 2).Value = 532
$ws.Cells.Item(77, 2).Value = 536
$ws.Cells.Item(78, 2).Value = 539
$ws.Cells.Item(79, 2).Value = 542
$ws.Cells.Item(80, 2).Value = 546
$ws.Cells.Item(81, 2).Value = 549
$ws.Cells.Item(82, 2).Value = 553
$ws.Cells.Item(83, 2).Value = 556
$ws.Cells.Item(84, 2).Value = 559
